$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ly86"
$ws.Range("C2").Value = "Cd180"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2604763333333334
$ws.Range("H2").Value = 0.7814290000000002
$ws.Range("I2").Value = 0.002013145958083725
$ws.Range("J2").Value = 0.002013145958083725
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.776312666666667
$ws.Range("N2").Value = 8.328938000000001
$ws.Range("O2").Value = 0.04197717588854396
$ws.Range("P2").Value = 0.04197717588854395
$ws.Range("Q2").Value = 0.7231637436002224
$ws.Range("R2").Value = 6.508473692402002
$ws.Range("S2").Value = 0.00008450618197179186
$ws.Range("T2").Value = 0.00008450618197179185

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ly86"
$ws.Range("C3").Value = "Cd180"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2604763333333334
$ws.Range("H3").Value = 0.7814290000000002
$ws.Range("I3").Value = 0.002013145958083725
$ws.Range("J3").Value = 0.002013145958083725
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1430836666666667
$ws.Range("N3").Value = 0.429251
$ws.Range("O3").Value = 0.002163390425926256
$ws.Range("P3").Value = 0.002163390425926256
$ws.Range("Q3").Value = 0.03726990885322223
$ws.Range("R3").Value = 0.335429179679
$ws.Range("S3").Value = 0.000004355220691710471
$ws.Range("T3").Value = 0.00000435522069171047

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ly86"
$ws.Range("C4").Value = "Cd180"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2604763333333334
$ws.Range("H4").Value = 0.7814290000000002
$ws.Range("I4").Value = 0.002013145958083725
$ws.Range("J4").Value = 0.002013145958083725
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 63.21922800000001
$ws.Range("N4").Value = 189.657684
$ws.Range("O4").Value = 0.9558594336855298
$ws.Range("P4").Value = 0.9558594336855297
$ws.Range("Q4").Value = 16.46711270560401
$ws.Range("R4").Value = 148.204014350436
$ws.Range("S4").Value = 0.001924284555420223
$ws.Range("T4").Value = 0.001924284555420222

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ly86"
$ws.Range("C5").Value = "Cd180"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4933343333333333
$ws.Range("H5").Value = 1.480003
$ws.Range("I5").Value = 0.003812837836069287
$ws.Range("J5").Value = 0.003812837836069287
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.776312666666667
$ws.Range("N5").Value = 8.328938000000001
$ws.Range("O5").Value = 0.04197717588854396
$ws.Range("P5").Value = 0.04197717588854395
$ws.Range("Q5").Value = 1.369650358534889
$ws.Range("R5").Value = 12.326853226814
$ws.Range("S5").Value = 0.0001600521644791758
$ws.Range("T5").Value = 0.0001600521644791758

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ly86"
$ws.Range("C6").Value = "Cd180"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4933343333333333
$ws.Range("H6").Value = 1.480003
$ws.Range("I6").Value = 0.003812837836069287
$ws.Range("J6").Value = 0.003812837836069287
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1430836666666667
$ws.Range("N6").Value = 0.429251
$ws.Range("O6").Value = 0.002163390425926256
$ws.Range("P6").Value = 0.002163390425926256
$ws.Range("Q6").Value = 0.07058808530588889
$ws.Range("R6").Value = 0.6352927677529999
$ws.Range("S6").Value = 0.000008248656870161679
$ws.Range("T6").Value = 0.000008248656870161677

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ly86"
$ws.Range("C7").Value = "Cd180"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4933343333333333
$ws.Range("H7").Value = 1.480003
$ws.Range("I7").Value = 0.003812837836069287
$ws.Range("J7").Value = 0.003812837836069287
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 63.21922800000001
$ws.Range("N7").Value = 189.657684
$ws.Range("O7").Value = 0.9558594336855298
$ws.Range("P7").Value = 0.9558594336855297
$ws.Range("Q7").Value = 31.188215699228
$ws.Range("R7").Value = 280.693941293052
$ws.Range("S7").Value = 0.003644537014719949
$ws.Range("T7").Value = 0.003644537014719949

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Ly86"
$ws.Range("C8").Value = "Cd180"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 128.6338933333333
$ws.Range("H8").Value = 385.9016799999999
$ws.Range("I8").Value = 0.9941740162058469
$ws.Range("J8").Value = 0.994174016205847
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.776312666666667
$ws.Range("N8").Value = 8.328938000000001
$ws.Range("O8").Value = 0.04197717588854396
$ws.Range("P8").Value = 0.04197717588854395
$ws.Range("Q8").Value = 357.1279074239822
$ws.Range("R8").Value = 3214.15116681584
$ws.Range("S8").Value = 0.04173261754209299
$ws.Range("T8").Value = 0.04173261754209299

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Ly86"
$ws.Range("C9").Value = "Cd180"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 128.6338933333333
$ws.Range("H9").Value = 385.9016799999999
$ws.Range("I9").Value = 0.9941740162058469
$ws.Range("J9").Value = 0.994174016205847
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1430836666666667
$ws.Range("N9").Value = 0.429251
$ws.Range("O9").Value = 0.002163390425926256
$ws.Range("P9").Value = 0.002163390425926256
$ws.Range("Q9").Value = 18.40540911574222
$ws.Range("R9").Value = 165.64868204168
$ws.Range("S9").Value = 0.002150786548364384
$ws.Range("T9").Value = 0.002150786548364384

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Ly86"
$ws.Range("C10").Value = "Cd180"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 128.6338933333333
$ws.Range("H10").Value = 385.9016799999999
$ws.Range("I10").Value = 0.9941740162058469
$ws.Range("J10").Value = 0.994174016205847
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 63.21922800000001
$ws.Range("N10").Value = 189.657684
$ws.Range("O10").Value = 0.9558594336855298
$ws.Range("P10").Value = 0.9558594336855297
$ws.Range("Q10").Value = 8132.13543116768
$ws.Range("R10").Value = 73189.21888050911
$ws.Range("S10").Value = 0.9502906121153896
$ws.Range("T10").Value = 0.9502906121153896
